$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"
#
# Column G is header "K" (row 1). Its per-row values were recomputed from the
# underlying strike data (s_vals) and are written back as static numbers,
# replacing the previous "Strike#"-derived figures. Map of row -> new K value
# (rows not listed, e.g. 30 and 70, already held the correct recalculated
# value of 0 and are left untouched).
$kValues = @{
    2 = 2;  3 = 2;  4 = 1;  5 = 1;  6 = 1;  7 = 2;  8 = 0;  9 = 2;  10 = 1;
    11 = 1; 12 = 2; 13 = 1; 14 = 0; 15 = 0; 16 = 1; 17 = 1; 18 = 1; 19 = 2;
    20 = 2; 21 = 1; 22 = 2; 23 = 2; 24 = 0; 25 = 0; 26 = 0; 27 = 2; 28 = 0;
    29 = 0; 31 = 1; 32 = 0; 33 = 2; 34 = 1; 35 = 0; 36 = 0; 37 = 0; 38 = 0;
    39 = 1; 40 = 1; 41 = 1; 42 = 1; 43 = 1; 44 = 1; 45 = 0; 46 = 1; 47 = 2;
    48 = 1; 49 = 1; 50 = 1; 51 = 2; 52 = 0; 53 = 1; 54 = 1; 55 = 0; 56 = 3;
    57 = 0; 58 = 1; 59 = 1; 60 = 0; 61 = 2; 62 = 1; 63 = 0; 64 = 1; 65 = 0;
    66 = 0; 67 = 1; 68 = 3; 69 = 1; 71 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
